# Auto-generated Excel COM-interop script applying the diff
# Updates market-data columns (H-N) across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 149.5
$ws.Range("J9").Value = 200
$ws.Range("L9").Value = 200
$ws.Range("N9").Value = -538

$ws.Range("H21").Value = 10019
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H23").Value = 10019
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("H34").Value = 401.16666
$ws.Range("I34").Value = 471.6
$ws.Range("K34").Value = 471.6
$ws.Range("M34").Value = -268.6

$ws.Range("H36").Value = 401.16666
$ws.Range("I36").Value = 471.6
$ws.Range("K36").Value = 471.6
$ws.Range("M36").Value = 243.4

$ws.Range("H62").Value = 1000
$ws.Range("I62").Value = 1000
$ws.Range("K62").Value = 1000
$ws.Range("M62").Value = -376

$ws.Range("H65").Value = 1000
$ws.Range("I65").Value = 1000
$ws.Range("K65").Value = 5000
$ws.Range("M65").Value = -1880

$ws.Range("H135").Value = 1062.5
$ws.Range("I135").Value = 775
$ws.Range("J135").Value = 1350
$ws.Range("K135").Value = 6975
$ws.Range("L135").Value = 12150
$ws.Range("M135").Value = -4440
$ws.Range("N135").Value = -17220

$ws.Range("H137").Value = 6673
$ws.Range("I137").Value = 3914.2856
$ws.Range("J137").Value = 11500.75
$ws.Range("K137").Value = 11742.8568
$ws.Range("L137").Value = 34502.25
$ws.Range("M137").Value = -9192.856800000001
$ws.Range("N137").Value = -39602.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 1250
$ws.Range("I21").Value = 1000
$ws.Range("J21").Value = 1500
$ws.Range("K21").Value = 1000
$ws.Range("L21").Value = 1500
$ws.Range("M21").Value = -626
$ws.Range("N21").Value = -2248

$ws.Range("H26").Value = 500
$ws.Range("I26").Value = 500
$ws.Range("K26").Value = 500
$ws.Range("M26").Value = -170

$ws.Range("H32").Value = 7113
$ws.Range("I32").Value = 2380.8
$ws.Range("K32").Value = 2380.8
$ws.Range("M32").Value = -2093.8

$ws.Range("H61").Value = 6769.8
$ws.Range("I61").Value = 3466.3333
$ws.Range("K61").Value = 3466.3333
$ws.Range("M61").Value = -3254.3333

$ws.Range("H101").Value = 20000
$ws.Range("J101").Value = 20000
$ws.Range("L101").Value = 20000
$ws.Range("N101").Value = -26490

$ws.Range("H136").Value = 6769.8
$ws.Range("I136").Value = 3466.3333
$ws.Range("K136").Value = 10398.9999
$ws.Range("M136").Value = -7848.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 45143.332
$ws.Range("I96").Value = 35428
$ws.Range("J96").Value = 50001
$ws.Range("K96").Value = 35428
$ws.Range("L96").Value = 50001
$ws.Range("M96").Value = -32682
$ws.Range("N96").Value = -55493

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 13999.333
$ws.Range("J58").Value = 16999
$ws.Range("L58").Value = 16999
$ws.Range("N58").Value = -17405

$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H107").Value = 620.8333
$ws.Range("I107").Value = 695.44446
$ws.Range("J107").Value = 397
$ws.Range("K107").Value = 695.44446
$ws.Range("L107").Value = 397
$ws.Range("M107").Value = 1224.55554
$ws.Range("N107").Value = -4237

$ws.Range("H136").Value = 13999.333
$ws.Range("J136").Value = 16999
$ws.Range("L136").Value = 50997
$ws.Range("N136").Value = -56097

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 100
$ws.Range("J11").Value = 100
$ws.Range("L11").Value = 300
$ws.Range("N11").Value = -580

$ws.Range("H17").Value = 795.1667
$ws.Range("J17").Value = 1490
$ws.Range("L17").Value = 4470
$ws.Range("N17").Value = -4808

$ws.Range("H26").Value = 1000
$ws.Range("J26").Value = 1000
$ws.Range("L26").Value = 3000
$ws.Range("N26").Value = -3576

$ws.Range("H39").Value = 7333.3335
$ws.Range("J39").Value = 7333.3335
$ws.Range("L39").Value = 22000.0005
$ws.Range("N39").Value = -22588.0005

$ws.Range("H55").Value = 2072.1428
$ws.Range("J55").Value = 2072.1428
$ws.Range("L55").Value = 6216.428400000001
$ws.Range("N55").Value = -6570.428400000001

$ws.Range("H58").Value = 5499.5
$ws.Range("I58").Value = 1000
$ws.Range("K58").Value = 3000
$ws.Range("M58").Value = -2872

$ws.Range("H68").Value = 1201.5

$ws.Range("H71").Value = 1201.5

$ws.Range("H129").Value = 1099.75
$ws.Range("I129").Value = 750
$ws.Range("K129").Value = 2250
$ws.Range("M129").Value = 2750

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1143800
$ws.Range("J11").Value = 1000320
$ws.Range("L11").Value = 1000320
$ws.Range("N11").Value = -1000598

$ws.Range("H101").Value = 32999
$ws.Range("J101").Value = 32999
$ws.Range("L101").Value = 32999
$ws.Range("N101").Value = -39489

$ws.Range("H126").Value = 1996
$ws.Range("I126").Value = 1996
$ws.Range("K126").Value = 5988
$ws.Range("M126").Value = -3518

$ws.Range("H132").Value = 9610.071
$ws.Range("I132").Value = 6755.222
$ws.Range("J132").Value = 14748.8
$ws.Range("K132").Value = 20265.666
$ws.Range("L132").Value = 44246.39999999999
$ws.Range("M132").Value = -17735.666
$ws.Range("N132").Value = -49306.39999999999

$ws.Range("H138").Value = 90000
$ws.Range("J138").Value = 90000
$ws.Range("L138").Value = 90000
$ws.Range("N138").Value = -100280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -205
$ws.Range("N22").Value = -1090

$ws.Range("H27").Value = 500
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 500
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 500
$ws.Range("M27").Value = -393
$ws.Range("N27").Value = -714

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 21666.334
$ws.Range("I75").Value = 14999
$ws.Range("K75").Value = 14999
$ws.Range("M75").Value = -14063

$ws.Range("H78").Value = 21666.334
$ws.Range("I78").Value = 14999
$ws.Range("K78").Value = 44997
$ws.Range("M78").Value = -40317

$ws.Range("H100").Value = 2133.3333
$ws.Range("I100").Value = 1500
$ws.Range("K100").Value = 3000
$ws.Range("M100").Value = -2459
